# [Fonds de solidarite] Add 2020-08-21 data
# Updates "nombre_aides" (column C) and "montant_total" (column D) figures
# for the rows whose underlying counts grew with the new data pull.
#
# The source cells are stored as text (numeric-looking strings), so we
# force the NumberFormat to Text ("@") before writing the new value -
# otherwise Excel/the COM layer would auto-coerce the assignment into a
# real number (losing the original "123.40"-style trailing-zero
# formatting to floating point rounding).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  C = "675";  D = "1551687.79" },
    @{ Row = 4;  C = "1015"; D = "3579453.47" },
    @{ Row = 6;  C = "651";  D = "2097277.78" },
    @{ Row = 9;  C = "182";  D = "545891.60" },
    @{ Row = 10; C = "364";  D = "1357808.18" },
    @{ Row = 11; C = "169";  D = "613307.11" },
    @{ Row = 12; C = "7";    D = "17950.00" },
    @{ Row = 14; C = "220";  D = "588362.00" },
    @{ Row = 16; C = "496";  D = "1826574.13" },
    @{ Row = 17; C = "142";  D = "440500.81" },
    @{ Row = 28; C = "280";  D = "724137.45" },
    @{ Row = 30; C = "558";  D = "2302182.70" },
    @{ Row = 32; C = "388";  D = "1328198.57" },
    @{ Row = 34; C = "16";   D = "47932.00" },
    @{ Row = 45; C = "378";  D = "1023772.74" },
    @{ Row = 47; C = "616";  D = "2408660.99" },
    @{ Row = 48; C = "417";  D = "1428367.16" },
    @{ Row = 50; C = "27";   D = "143011.07" },
    @{ Row = 76; C = "926";  D = "3228740.26" }
)

foreach ($u in $updates) {
    $cCell = $ws.Cells.Item($u.Row, 3)
    $dCell = $ws.Cells.Item($u.Row, 4)

    $cCell.NumberFormat = "@"
    $dCell.NumberFormat = "@"

    $cCell.Value = $u.C
    $dCell.Value = $u.D
}
